$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.738254
$ws.Range("H2").Value = 2.214762
$ws.Range("I2").Value = 0.005691320045803731
$ws.Range("J2").Value = 0.005691320045803731
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 0.2502247950986666
$ws.Range("R2").Value = 2.252023155888
$ws.Range("S2").Value = 0.005691320045803731
$ws.Range("T2").Value = 0.005691320045803731

# Row 3
$ws.Range("I3").Value = 0.9440493064670392
$ws.Range("J3").Value = 0.9440493064670391
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 41.506107963112
$ws.Range("R3").Value = 373.554971668008
$ws.Range("S3").Value = 0.9440493064670392
$ws.Range("T3").Value = 0.9440493064670391

# Row 4
$ws.Range("G4").Value = 5.698467
$ws.Range("H4").Value = 17.095401
$ws.Range("I4").Value = 0.0439304080539368
$ws.Range("J4").Value = 0.04393040805393679
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3389413333333333
$ws.Range("N4").Value = 1.016824
$ws.Range("Q4").Value = 1.931446002936
$ws.Range("R4").Value = 17.383014026424
$ws.Range("S4").Value = 0.0439304080539368
$ws.Range("T4").Value = 0.04393040805393679

# Row 5
$ws.Range("G5").Value = 0.8209666666666666
$ws.Range("H5").Value = 2.4629
$ws.Range("I5").Value = 0.006328965433220369
$ws.Range("J5").Value = 0.006328965433220369
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3389413333333333
$ws.Range("N5").Value = 1.016824
$ws.Range("Q5").Value = 0.2782595366222222
$ws.Range("R5").Value = 2.5043358296
$ws.Range("S5").Value = 0.006328965433220369
$ws.Range("T5").Value = 0.006328965433220369
